# Applies the cryptos.xlsx price/volume refresh described in the commit diff.
# Numeric-looking strings (e.g. "564.61") are written via a Text-formatted
# round-trip so Excel keeps them as literal strings (matching the inline-string
# cells openpyxl originally wrote) instead of silently parsing them into numbers;
# the cell style is then reset to "Normal" so no stray formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2: Bitcoin
$ws.Range("D2").Value = "62.284.48"
$ws.Range("E2").Value = "  +1.45%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "2.423.25"
$ws.Range("E3").Value = "  +1.85%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.08%  "

# Row 5: BNB
Set-TextValue "D5" "564.61"
$ws.Range("E5").Value = "  +2.32%  "

# Row 6: Solana
Set-TextValue "D6" "144.75"
$ws.Range("E6").Value = "  +3.53%  "

# Row 7: USDC
$ws.Range("E7").Value = "  -0.01%  "

# Row 8: XRP
Set-TextValue "D8" "0.533"
$ws.Range("E8").Value = "  +1.93%  "

# Row 9: LidoStakedEther
$ws.Range("D9").Value = "2.421.33"
$ws.Range("E9").Value = "  +1.72%  "

# Row 10: Dogecoin
$ws.Range("E10").Value = "  +1.98%  "

# Row 11: TRON
$ws.Range("E11").Value = "  -1.52%  "

# Row 12: Toncoin
Set-TextValue "D12" "5.39"
$ws.Range("E12").Value = "  +0.89%  "

# Row 13: Cardano
Set-TextValue "D13" "0.354"
$ws.Range("E13").Value = "  +0.73%  "

# Row 14: Avalanche
$ws.Range("E14").Value = "  +2.11%  "

# Row 15: ShibaInu
Set-TextValue "D15" "0.0000178"
$ws.Range("E15").Value = "  +6.09%  "

# Row 16: WrappedliquidstakedEther2.0
$ws.Range("D16").Value = "2.859.92"
$ws.Range("E16").Value = "  +1.84%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "62.003.58"
$ws.Range("E17").Value = "  +0.97%  "

# Row 18: WrappedEther
$ws.Range("D18").Value = "2.419.24"
$ws.Range("E18").Value = "  +1.53%  "

# Row 19: Chainlink
Set-TextValue "D19" "11.32"
$ws.Range("E19").Value = "  +3.00%  "

# Row 20: Polkadot
$ws.Range("E20").Value = "  +1.27%  "

# Row 21: BitcoinCash
Set-TextValue "D21" "325.03"
$ws.Range("E21").Value = "  +1.15%  "

# Row 22: Uniswap
Set-TextValue "D22" "6.76"
$ws.Range("E22").Value = "  +0.54%  "

# Row 23: Dai
$ws.Range("E23").Value = "  +0.05%  "

# Row 24: Litecoin
Set-TextValue "D24" "65.61"
$ws.Range("E24").Value = "  +1.80%  "

# Row 25: SuiNetwork
$ws.Range("E25").Value = "  -2.99%  "

# Row 26: Aptos
Set-TextValue "D26" "9.01"
$ws.Range("E26").Value = "  +0.95%  "

# Row 27: Bittensor
Set-TextValue "D27" "590.81"
$ws.Range("E27").Value = "  +14.32%  "

# Row 28: PEPE
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0950"
$ws.Range("E28").Value = "  +5.46%  "

# Row 29: Binance-PegBSC-USD
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D29" "1.00"
$ws.Range("E29").Value = "  +0.09%  "

# Row 30: WrappedeETH
$ws.Range("B30").Value = "WrappedeETH"
$ws.Range("C30").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D30").Value = "2.528.43"
$ws.Range("E30").Value = "  +1.28%  "

# Row 31: Fetch.AI
Set-TextValue "D31" "1.48"
$ws.Range("E31").Value = "  +5.82%  "

# Row 32: InternetComputer(DFINITY)
$ws.Range("E32").Value = "  +1.00%  "

# Row 33: Kaspa
$ws.Range("E33").Value = "  +0.69%  "

# Row 34: PancakeSwap
$ws.Range("E34").Value = "  +2.13%  "

# Row 35: ImmutableX
$ws.Range("E35").Value = "  +0.99%  "

# Row 36: RenderToken
Set-TextValue "D36" "5.76"
$ws.Range("E36").Value = "  +4.12%  "

# Row 37: FirstDigitalUSD
Set-TextValue "D37" "0.999"
$ws.Range("E37").Value = "  +0.04%  "

# Row 38: NEARProtocol
$ws.Range("E38").Value = "  +2.20%  "

# Row 39: Monero
Set-TextValue "D39" "153.73"
$ws.Range("E39").Value = "  +5.34%  "

# Row 40: PolygonEcosystemToken
$ws.Range("E40").Value = "  +1.30%  "

# Row 41: EthereumClassic
Set-TextValue "D41" "18.72"
$ws.Range("E41").Value = "  +1.03%  "

# Row 42: Stacks
Set-TextValue "D42" "1.84"
$ws.Range("E42").Value = "  -2.46%  "

# Row 43: USDe
$ws.Range("E43").Value = "  -0.15%  "

# Row 44: dogwifhat
$ws.Range("E44").Value = "  +9.16%  "

# Row 45: Aave
Set-TextValue "D45" "150.27"
$ws.Range("E45").Value = "  +1.40%  "

# Row 46: Filecoin
Set-TextValue "D46" "3.65"
$ws.Range("E46").Value = "  +1.44%  "

# Row 47: Hedera
$ws.Range("E47").Value = "  +2.35%  "

# Row 48: InjectiveProtocol
Set-TextValue "D48" "20.46"
$ws.Range("E48").Value = "  +4.03%  "

# Row 49: Mantle
$ws.Range("E49").Value = "  +1.99%  "

# Row 50: Stellar
$ws.Range("E50").Value = "  +2.23%  "

# Row 51: VeChain
$ws.Range("E51").Value = "  +1.88%  "
